$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("B2").Value = 0.000000000000000000145568428466758
$ws.Range("C2").Value = 0.0004851533198767937
$ws.Range("D2").Value = 0.1474043695038804
$ws.Range("E2").Value = 0.8131214864200552
$ws.Range("B3").Value = 0.000000001573017879546552
$ws.Range("C3").Value = 0.008234292916330184
$ws.Range("D3").Value = 0.1367730723658386
$ws.Range("E3").Value = 0.754982290267857
$ws.Range("B4").Value = 0.00000000002455525622174766
$ws.Range("C4").Value = 0.007637918818402064
$ws.Range("D4").Value = 0.1117744852483172
$ws.Range("E4").Value = 0.6070804060303006
$ws.Range("C5").Value = 0.00001687926027758808
$ws.Range("D5").Value = 0.007208659182546
$ws.Range("E5").Value = 0.04501130247973132

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C2").Value = 0.0004851533198768265
$ws.Range("D2").Value = 0.1123080910505831
$ws.Range("E2").Value = 0.8131214864201096
$ws.Range("C3").Value = 0.008234292916330738
$ws.Range("D3").Value = 0.1042080551358841
$ws.Range("E3").Value = 0.7549822902679078
$ws.Range("C4").Value = 0.007637918818402577
$ws.Range("D4").Value = 0.08516151257015217
$ws.Range("E4").Value = 0.6070804060303412
$ws.Range("C5").Value = 0.00001687926027758942
$ws.Range("D5").Value = 0.005492311758130723
$ws.Range("E5").Value = 0.04501130247973491

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = 0.0004996343988088803
$ws.Range("C2").Value = 0.04937646078717243
$ws.Range("D2").Value = 1.264699885066036
$ws.Range("E2").Value = 2.070245676915336
$ws.Range("B3").Value = 0.003689272420938101
$ws.Range("C3").Value = 0.05729477318728819
$ws.Range("D3").Value = 0.8886328477300586
$ws.Range("E3").Value = 1.552184907137255
$ws.Range("B4").Value = 0.009318027231660879
$ws.Range("C4").Value = 0.0440480498952964
$ws.Range("D4").Value = 1.05784163086344
$ws.Range("E4").Value = 1.788123587539221
$ws.Range("B5").Value = 0.003016841357489923
$ws.Range("C5").Value = 0.04146064243109764
$ws.Range("D5").Value = 0.9969023953454611
$ws.Range("E5").Value = 1.757829693982964

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = 0.00008341175675907166
$ws.Range("C2").Value = 0.01652097409187206
$ws.Range("D2").Value = 0.836659011110328
$ws.Range("E2").Value = 2.553279458085891
$ws.Range("B3").Value = 0.00009666419205255739
$ws.Range("C3").Value = 0.01561751262125499
$ws.Range("D3").Value = 0.4567345111362074
$ws.Range("E3").Value = 1.536080615865836
$ws.Range("B4").Value = 0.0005273422179231309
$ws.Range("C4").Value = 0.01397792516947091
$ws.Range("D4").Value = 0.6174460803910278
$ws.Range("E4").Value = 2.128650746148338
$ws.Range("B5").Value = 0.0002918342513006147
$ws.Range("C5").Value = 0.0161799119914323
$ws.Range("D5").Value = 0.7374355342131388
$ws.Range("E5").Value = 2.277846581435195

